$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 387.1016018289494
$ws.Range("D2").Value = 387.1268169634347
$ws.Range("C3").Value = 387.5687990774547
$ws.Range("D3").Value = 387.5917350937759
$ws.Range("C4").Value = 387.1826566112436
$ws.Range("D4").Value = 387.2390567885702
$ws.Range("C5").Value = 386.9414880307917
$ws.Range("D5").Value = 386.9838668866848
$ws.Range("C6").Value = 387.2625050464738
$ws.Range("D6").Value = 387.3690273265964
$ws.Range("C7").Value = 386.8173255130815
$ws.Range("D7").Value = 386.8964301981485
$ws.Range("C8").Value = 387.439487283573
$ws.Range("D8").Value = 387.6035205697006
$ws.Range("C9").Value = 387.4388173465784
$ws.Range("D9").Value = 387.5365936072684
$ws.Range("C10").Value = 386.9694360693297
$ws.Range("D10").Value = 386.9888011663205
$ws.Range("C11").Value = 387.1909790238212
$ws.Range("D11").Value = 387.3076946489816
$ws.Range("E11").Value = 1.192266613200502
$ws.Range("F11").Value = -0.6937070539837423
$ws.Range("G11").Value = 1.380739891101147
$ws.Range("C12").Value = 387.5317655884826
$ws.Range("D12").Value = 387.6479332138931
$ws.Range("E12").Value = 1.333573253186748
$ws.Range("F12").Value = -0.2459064695264767
$ws.Range("G12").Value = 0.2357726206422357
$ws.Range("C13").Value = 387.345868595697
$ws.Range("D13").Value = 387.4181034194296
$ws.Range("E13").Value = 1.463258557030148
$ws.Range("F13").Value = -0.4569013605566123
$ws.Range("G13").Value = 0.9812587943120659
$ws.Range("C14").Value = 387.1412462094812
$ws.Range("D14").Value = 387.217081261914
$ws.Range("E14").Value = 1.432882181232127
$ws.Range("F14").Value = -0.4359629151490947
$ws.Range("G14").Value = 0.7483782737744353
$ws.Range("C15").Value = 388.4381883895534
$ws.Range("D15").Value = 388.4803554515964
$ws.Range("E15").Value = 1.445688955703633
$ws.Range("F15").Value = -0.5779276682639747
$ws.Range("G15").Value = 2.11375978505539
$ws.Range("C16").Value = 388.4124658794146
$ws.Range("D16").Value = 388.479071799711
$ws.Range("E16").Value = 1.587305773226825
$ws.Range("F16").Value = -0.1581859639359714
$ws.Range("G16").Value = 0.1198330976075468
$ws.Range("C17").Value = 387.7634293206705
$ws.Range("D17").Value = 387.6771289681381
$ws.Range("E17").Value = 1.69418452977483
$ws.Range("F17").Value = 0.3939887492320298
$ws.Range("G17").Value = 0.152294289086699
$ws.Range("C18").Value = 387.1302126342749
$ws.Range("D18").Value = 387.1716338904715
$ws.Range("E18").Value = 1.351063952638233
$ws.Range("F18").Value = -0.1214769430852471
$ws.Range("G18").Value = 2.089503294318475
$ws.Range("C19").Value = 387.9329004408691
$ws.Range("D19").Value = 387.8003093450315
